$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("LuSTRE")
$ws2 = $wb.Worksheets.Item("Getty ULAN")
$ws3 = $wb.Worksheets.Item("Getty AAT ")
$ws4 = $wb.Worksheets.Item("Getty TGN")

# The "name" column (A) data value on each Getty sheet duplicates the sheet's
# own display name -- update it to the whitespace-free spelling too.
$ws2.Range("A3").Value = "GettyULAN"
$ws3.Range("A3").Value = "GettyAAT"
$ws4.Range("A3").Value = "GettyTGN"

# Remove whitespace from the sheet tab names (AAT keeps its trailing space).
$ws2.Name = "GettyULAN"
$ws3.Name = "GettyAAT "
$ws4.Name = "GettyTGN"

# Move the remembered selection on each Getty sheet.
$ws2.Activate()
$ws2.Range("A3").Select()

$ws3.Activate()
$ws3.Range("B16").Select()

$ws4.Activate()
$ws4.Range("E28").Select()

# GettyTGN ends up the active / visible tab when the workbook is saved.
$ws4.Activate()
